# Atualização de bases das ligas, do dia: 17-06-2024 às 21:10
# Swap the data of several pairs of match rows (everything except the
# row's "id" label in column A, and the Div/Date columns C/D, which are
# identical for both rows of each pair anyway).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that must be swapped between the two rows of each pair:
# B (match id), E..AD (HomeTeam .. PL_AhUnder)
$cols = @(2) + (5..30)

$rowPairs = @(
    @(39, 40),
    @(60, 61),
    @(77, 78),
    @(88, 89),
    @(186, 187),
    @(260, 261)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($c in $cols) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)
        $v1 = $cell1.Value()
        $v2 = $cell2.Value()
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}
